$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 11 with new activity entry
$ws.Range("A11").Value = "Di"
$ws.Range("B11").Value = "10.12.2019"
$ws.Range("C11").Value = 0.32291666666666669
$ws.Range("D11").Value = 0.54166666666666663
$ws.Range("E11").Value = "Sprintplaning und Programierung"

# Update selection to E15 as recorded in the sheet view
$ws.Range("E15").Select()
